$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 and 3 (2008年 and 2009年), shifting 2010年/2011年 up.
$ws.Rows("2:3").Delete()
